# Applies the "error handling for missing sheets in xlsx file" edit:
#  - Rename "Descriptive Metadata" -> "Metadata"
#  - Remove the row of real data that used to live at row 35 of the Metadata sheet
#  - Delete the "Type Of Resource Values" and "Source and Location Values" lookup sheets
#  - Point the defined names that used to target those sheets at #REF! (their source is gone)
#  - Add a new "Xsls" sheet (after "Mappings") listing the extra mapping keys
#  - Make "Xsls" the active tab

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- rename the main data-entry sheet ---
$metadata = $wb.Worksheets.Item("Descriptive Metadata")
$metadata.Name = "Metadata"

# --- the last data row (row 35) is removed from Metadata (shifts nothing up below it, it's the last row) ---
$metadata.Rows.Item(35).Delete()

# --- drop the now-orphaned lookup sheets ---
$wb.Worksheets.Item("Type Of Resource Values").Delete()
$wb.Worksheets.Item("Source and Location Values").Delete()

# --- the defined names that referenced those sheets now dangle ---
foreach ($n in $wb.Names) {
    $n.RefersTo = "=#REF!"
}

# --- add the new "Xsls" sheet right after "Mappings" ---
$mappings = $wb.Worksheets.Item("Mappings")
$xsls = $wb.Worksheets.Add($null, $mappings)
$xsls.Name = "Xsls"

$values = @("titleNonSort", "multiNamePart", "blankNamePart", "dateCreatedSplit", "subjectSplit", "normalizeDate", "OrderedTemplates", "blankNodes")
for ($i = 0; $i -lt $values.Length; $i++) {
    $xsls.Cells.Item($i + 1, 1).Value = $values[$i]
}

# --- Xsls becomes the active tab ---
$xsls.Activate()
